# Rename the inline logo pictures living in the document's headers/footers.
#
#   headers (BTec_Logo-Orange, .jpg):  image1.jpg -> image2.jpg
#   footers (PearsonLogo, .png):       image2.png -> image1.png
#
# The two headers/footers are "default" (primary) and "first page"
# (this document has Different First Page turned on), reached through
# Section 1's Headers/Footers collections:
#   Headers./Footers.Item(1) = wdHeaderFooterPrimary
#   Headers./Footers.Item(2) = wdHeaderFooterFirstPage

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-LogoInStory($story, $newName) {
    if ($story.Exists) {
        $shape = $story.Range.InlineShapes.Item(1)
        $shape.Name = $newName
    }
}

# Headers - BTec logo, image1.jpg -> image2.jpg
Rename-LogoInStory $sec.Headers.Item(1) "image2.jpg"
Rename-LogoInStory $sec.Headers.Item(2) "image2.jpg"

# Footers - Pearson logo, image2.png -> image1.png
Rename-LogoInStory $sec.Footers.Item(1) "image1.png"
Rename-LogoInStory $sec.Footers.Item(2) "image1.png"
